$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings; all runs share the same font, so a plain value set is safe) ---
$ws.Range("A8").Value = "Volume 31   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/25/2024  Through  3/31/2024"

# --- Weekly crime-stat numeric updates (rows 16-30) ---
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 44.444444444444
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = 14.634146341463
$ws.Range("L16").Value = 30.555555555555
$ws.Range("M16").Value = -22.950819672131
$ws.Range("N16").Value = -75.132275132275
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 78.260869565217
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 79
$ws.Range("K17").Value = 26.582278481012
$ws.Range("L17").Value = 61.290322580645
$ws.Range("M17").Value = 156.410256410256
$ws.Range("N17").Value = -29.577464788732
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 340
$ws.Range("I18").Value = 47
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = 62.068965517241
$ws.Range("L18").Value = 38.235294117647
$ws.Range("M18").Value = 176.470588235294
$ws.Range("N18").Value = -53.921568627451
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -32.432432432432
$ws.Range("I19").Value = 91
$ws.Range("J19").Value = 113
$ws.Range("K19").Value = -19.469026548672
$ws.Range("L19").Value = -7.142857142857
$ws.Range("M19").Value = 56.896551724137
$ws.Range("N19").Value = -44.171779141104
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -11.111111111111
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -42.857142857142
$ws.Range("L20").Value = -44.827586206896
$ws.Range("M20").Value = 220
$ws.Range("N20").Value = -82.608695652173
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 45
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = 30.952380952381
$ws.Range("I21").Value = 303
$ws.Range("J21").Value = 298
$ws.Range("K21").Value = 1.677852348993
$ws.Range("L21").Value = 12.639405204461
$ws.Range("M21").Value = 65.573770491803
$ws.Range("N21").Value = -57.383966244725
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 33.333333333333
$ws.Range("F23").Value = 27
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = 3.846153846153
$ws.Range("I23").Value = 85
$ws.Range("J23").Value = 85
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1.190476190476
$ws.Range("M23").Value = 73.469387755102
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -22.222222222222
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = -14.864864864864
$ws.Range("I24").Value = 188
$ws.Range("J24").Value = 203
$ws.Range("K24").Value = -7.389162561576
$ws.Range("L24").Value = 16.049382716049
$ws.Range("M24").Value = 37.226277372262
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 12
$ws.Range("H25").Value = -42.857142857142
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = -25.581395348837
$ws.Range("L25").Value = -3.030303030303
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 60
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 165
$ws.Range("J26").Value = 130
$ws.Range("K26").Value = 26.923076923076
$ws.Range("L26").Value = 32
$ws.Range("M26").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 66.666666666666
$ws.Range("L28").Value = -7.142857142857
$ws.Range("L29").Value = -83.333333333333
$ws.Range("N29").Value = -96.296296296296
$ws.Range("L30").Value = -75
$ws.Range("N30").Value = -96.153846153846

# --- Row 28: C28/D28 become the "0" N/A placeholder, E28 becomes the "***.*" N/A placeholder ---
# Reuse the format+value from C22/D22/E22, which already hold the identical N/A placeholders/style,
# so the shared-string + style table stay exactly as Excel would produce them natively.
$ws.Range("C22").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C22").Copy()
$ws.Paste($ws.Range("C28"))

$ws.Range("D22").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D22").Copy()
$ws.Paste($ws.Range("D28"))

$ws.Range("E22").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E22").Copy()
$ws.Paste($ws.Range("E28"))

$excel.CutCopyMode = 0
